$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the chequing-account test plan rows (E7:G14) with the revised
# --- Module 2 Assignment 2 test cases / inputs / expected results.

# Row 7 - Valid account data
$ws.Range("E7").Value = "Valid account data"
$ws.Range("F7").Value = "(`"1234567`", `"C001`", 1559.49, date(2024, 1, 1), -100.0, 0.05)"
$ws.Range("G7").Value = "account_number=`"1234567`", balance=1559.49, date_created=date(2024, 1, 1), overdraft_limit=-100.0, overdraft_rate=0.05"

# Row 8 - overdraft_limit invalid type
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "(`"1234567`", `"C001`", 1559.49, date(2024, 1, 1), `"invalid`", 0.05)"
$ws.Range("G8").Value = "overdraft_limit defaults to -100.0"

# Row 9 - overdraft_rate invalid type
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "(`"1234567`", `"C001`", 1559.49, date(2024, 1, 1), -100.0, `"invalid`")"
$ws.Range("G9").Value = "overdraft_rate defaults to 0.05"

# Row 10 - date_created invalid type
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = "(`"1234567`", `"C001`", 1559.49, `"invalid`", -100.0, 0.05)"
$ws.Range("G10").Value = "date_created defaults to today’s date (e.g., 2025-02-10)"

# Row 11 - get_service_charges(), balance greater than overdraft limit
$ws.Range("E11").Value = "self.balance = 500"
$ws.Range("F11").Value = "get_service_charges()"
$ws.Range("G11").Value = "Returns base service charge of 0.50"

# Row 12 - get_service_charges(), balance less than overdraft limit
$ws.Range("E12").Value = "self.balance = -300"
$ws.Range("F12").Value = "get_service_charges()"
$ws.Range("G12").Value = "Returns 15.50"

# Row 13 - get_service_charges(), balance equal to overdraft limit
$ws.Range("E13").Value = "self.balance = -100"
$ws.Range("F13").Value = "get_service_charges()"
$ws.Range("G13").Value = "Returns base service charge of 0.50"

# Row 14 - str(chequing_account)
$ws.Range("E14").Value = "Instance initialized with known values"
$ws.Range("F14").Value = "str(chequing_account)"
$ws.Range("G14").Value = "Returns `"Account Number: 1234567 Balance: `$1,559.49
Overdraft Limit: `$-100.00 Overdraft Rate: 5.00% Account Type: Cheq`""

# --- Update the saved selection / scroll position of the sheet view.
$ws.Range("E14").Select() | Out-Null
